# Auto-generated: update FFXIV Leve-profit market-data cells per scheduled runner diff.
# For each touched sheet, write the new numeric values into H:N columns for the affected rows.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2611.074
$ws.Range("I51").Value = 2321.524
$ws.Range("K51").Value = 2321.524
$ws.Range("M51").Value = -1837.524
$ws.Range("H64").Value = 7531.0713
$ws.Range("I64").Value = 3975
$ws.Range("J64").Value = 10198.125
$ws.Range("K64").Value = 3975
$ws.Range("L64").Value = 10198.125
$ws.Range("M64").Value = -3727
$ws.Range("N64").Value = -10694.125
$ws.Range("H67").Value = 7531.0713
$ws.Range("I67").Value = 3975
$ws.Range("J67").Value = 10198.125
$ws.Range("K67").Value = 3975
$ws.Range("L67").Value = 10198.125
$ws.Range("M67").Value = -3117
$ws.Range("N67").Value = -11914.125
$ws.Range("H70").Value = 70020.336
$ws.Range("I70").Value = 2877.6667
$ws.Range("J70").Value = 86806
$ws.Range("K70").Value = 8633.000100000001
$ws.Range("L70").Value = 260418
$ws.Range("M70").Value = -8363.000100000001
$ws.Range("N70").Value = -260958
$ws.Range("H73").Value = 70020.336
$ws.Range("I73").Value = 2877.6667
$ws.Range("J73").Value = 86806
$ws.Range("K73").Value = 8633.000100000001
$ws.Range("L73").Value = 260418
$ws.Range("M73").Value = -7697.000100000001
$ws.Range("N73").Value = -262290
$ws.Range("H74").Value = 153942.42
$ws.Range("I74").Value = 254399.25
$ws.Range("J74").Value = 20000
$ws.Range("K74").Value = 254399.25
$ws.Range("L74").Value = 20000
$ws.Range("M74").Value = -253463.25
$ws.Range("N74").Value = -21872
$ws.Range("H77").Value = 153942.42
$ws.Range("I77").Value = 254399.25
$ws.Range("J77").Value = 20000
$ws.Range("K77").Value = 1271996.25
$ws.Range("L77").Value = 100000
$ws.Range("M77").Value = -1267316.25
$ws.Range("N77").Value = -109360
$ws.Range("H80").Value = 306.9375
$ws.Range("I80").Value = 183.28572
$ws.Range("K80").Value = 549.85716
$ws.Range("M80").Value = 448.14284
$ws.Range("H83").Value = 306.9375
$ws.Range("I83").Value = 183.28572
$ws.Range("K83").Value = 1649.57148
$ws.Range("M83").Value = 3342.42852
$ws.Range("H86").Value = 2499.6667
$ws.Range("I86").Value = 1499
$ws.Range("K86").Value = 1499
$ws.Range("M86").Value = -376
$ws.Range("H89").Value = 2499.6667
$ws.Range("I89").Value = 1499
$ws.Range("K89").Value = 7495
$ws.Range("M89").Value = -1879
$ws.Range("H116").Value = 2854.2856
$ws.Range("I116").Value = 2768.3333
$ws.Range("J116").Value = 2918.75
$ws.Range("K116").Value = 2768.3333
$ws.Range("L116").Value = 2918.75
$ws.Range("M116").Value = 673.6667000000002
$ws.Range("N116").Value = -9802.75
$ws.Range("H132").Value = 1942.3704
$ws.Range("I132").Value = 1709.7084
$ws.Range("K132").Value = 5129.1252
$ws.Range("M132").Value = -2599.1252
$ws.Range("H138").Value = 6175370.5
$ws.Range("I138").Value = 1141.3478
$ws.Range("J138").Value = 8623772
$ws.Range("K138").Value = 3424.0434
$ws.Range("L138").Value = 25871316
$ws.Range("M138").Value = 1715.9566
$ws.Range("N138").Value = -25881596

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13600.893
$ws.Range("I32").Value = 9310.409
$ws.Range("K32").Value = 9310.409
$ws.Range("M32").Value = -9023.409

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2351.9167
$ws.Range("I86").Value = 2279.3667
$ws.Range("K86").Value = 2279.3667
$ws.Range("M86").Value = -1156.3667
$ws.Range("H89").Value = 2351.9167
$ws.Range("I89").Value = 2279.3667
$ws.Range("K89").Value = 11396.8335
$ws.Range("M89").Value = -5780.833500000001
$ws.Range("H99").Value = 10499.833
$ws.Range("I99").Value = 2510
$ws.Range("J99").Value = 12097.8
$ws.Range("K99").Value = 2510
$ws.Range("L99").Value = 12097.8
$ws.Range("M99").Value = -1012
$ws.Range("N99").Value = -15093.8

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 193.16667
$ws.Range("I22").Value = 183.45454
$ws.Range("K22").Value = 183.45454
$ws.Range("M22").Value = 166.54546
$ws.Range("H31").Value = 50587.76
$ws.Range("J31").Value = 3407.25
$ws.Range("L31").Value = 3407.25
$ws.Range("N31").Value = -3997.25
$ws.Range("H34").Value = 50587.76
$ws.Range("J34").Value = 3407.25
$ws.Range("L34").Value = 3407.25
$ws.Range("N34").Value = -3811.25
$ws.Range("H58").Value = 4171.1816
$ws.Range("I58").Value = 3999.75
$ws.Range("J58").Value = 4269.143
$ws.Range("K58").Value = 3999.75
$ws.Range("L58").Value = 4269.143
$ws.Range("M58").Value = -3796.75
$ws.Range("N58").Value = -4675.143
$ws.Range("H62").Value = 11525.167
$ws.Range("I62").Value = 12504
$ws.Range("J62").Value = 10154.8
$ws.Range("K62").Value = 12504
$ws.Range("L62").Value = 10154.8
$ws.Range("M62").Value = -11880
$ws.Range("N62").Value = -11402.8
$ws.Range("H65").Value = 11525.167
$ws.Range("I65").Value = 12504
$ws.Range("J65").Value = 10154.8
$ws.Range("K65").Value = 62520
$ws.Range("L65").Value = 50774
$ws.Range("M65").Value = -59400
$ws.Range("N65").Value = -57014
$ws.Range("H99").Value = 3671.3333
$ws.Range("I99").Value = 3000
$ws.Range("J99").Value = 5014
$ws.Range("K99").Value = 3000
$ws.Range("L99").Value = 5014
$ws.Range("M99").Value = -1502
$ws.Range("N99").Value = -8010
$ws.Range("H126").Value = 3671.3333
$ws.Range("I126").Value = 3000
$ws.Range("J126").Value = 5014
$ws.Range("K126").Value = 9000
$ws.Range("L126").Value = 15042
$ws.Range("M126").Value = -6530
$ws.Range("N126").Value = -19982
$ws.Range("H134").Value = 11317.293
$ws.Range("I134").Value = 5086.057
$ws.Range("K134").Value = 15258.171
$ws.Range("M134").Value = -12723.171
$ws.Range("H136").Value = 4171.1816
$ws.Range("I136").Value = 3999.75
$ws.Range("J136").Value = 4269.143
$ws.Range("K136").Value = 11999.25
$ws.Range("L136").Value = 12807.429
$ws.Range("M136").Value = -9449.25
$ws.Range("N136").Value = -17907.429

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 834.8889
$ws.Range("J113").Value = 868
$ws.Range("L113").Value = 2604
$ws.Range("N113").Value = -6944
$ws.Range("H137").Value = 2480.2942
$ws.Range("J137").Value = 4175.1665
$ws.Range("L137").Value = 12525.4995
$ws.Range("N137").Value = -22725.4995
$ws.Range("H139").Value = 2122
$ws.Range("I139").Value = 1829.3334
$ws.Range("K139").Value = 5488.0002
$ws.Range("M139").Value = -348.0002000000004

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3107
$ws.Range("I80").Value = 2598.1428
$ws.Range("K80").Value = 2598.1428
$ws.Range("M80").Value = -1600.1428
$ws.Range("H83").Value = 3107
$ws.Range("I83").Value = 2598.1428
$ws.Range("K83").Value = 12990.714
$ws.Range("M83").Value = -7998.714
$ws.Range("H97").Value = 705.2105
$ws.Range("I97").Value = 851.93335
$ws.Range("J97").Value = 155
$ws.Range("K97").Value = 851.93335
$ws.Range("L97").Value = 155
$ws.Range("M97").Value = -355.93335
$ws.Range("N97").Value = -1147
$ws.Range("H102").Value = 47621904
$ws.Range("I102").Value = 2831.2222
$ws.Range("J102").Value = 333336320
$ws.Range("K102").Value = 2831.2222
$ws.Range("L102").Value = 333336320
$ws.Range("M102").Value = -1209.2222
$ws.Range("N102").Value = -333339564

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4006.64
$ws.Range("I136").Value = 3528.9565
$ws.Range("J136").Value = 9500
$ws.Range("K136").Value = 10586.8695
$ws.Range("L136").Value = 28500
$ws.Range("M136").Value = -8036.869499999999
$ws.Range("N136").Value = -33600
$ws.Range("H141").Value = 109993
$ws.Range("J141").Value = 109993
$ws.Range("L141").Value = 109993
$ws.Range("N141").Value = -120353

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 57809.5
$ws.Range("J46").Value = 57809.5
$ws.Range("L46").Value = 57809.5
$ws.Range("N46").Value = -58271.5
$ws.Range("H122").Value = 1888.5483
$ws.Range("I122").Value = 1782.35
$ws.Range("J122").Value = 2081.6365
$ws.Range("K122").Value = 5347.049999999999
$ws.Range("L122").Value = 6244.9095
$ws.Range("M122").Value = -2897.049999999999
$ws.Range("N122").Value = -11144.9095
$ws.Range("H126").Value = 12445.637
$ws.Range("I126").Value = 14988.5
$ws.Range("J126").Value = 5664.6665
$ws.Range("K126").Value = 44965.5
$ws.Range("L126").Value = 16993.9995
$ws.Range("M126").Value = -42495.5
$ws.Range("N126").Value = -21933.9995
$ws.Range("H134").Value = 57809.5
$ws.Range("J134").Value = 57809.5
$ws.Range("L134").Value = 173428.5
$ws.Range("N134").Value = -178498.5
$ws.Range("H138").Value = 55166.668
$ws.Range("J138").Value = 55166.668
$ws.Range("L138").Value = 55166.668
$ws.Range("N138").Value = -65446.668

